$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "L1cam"
$ws.Cells.Item(2,3).Value = "Erbb2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 11.88712433333333
$ws.Cells.Item(2,8).Value = 35.661373
$ws.Cells.Item(2,9).Value = 0.5967229292030898
$ws.Cells.Item(2,10).Value = 0.5967229292030898
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 3.155977333333333
$ws.Cells.Item(2,14).Value = 9.467931999999999
$ws.Cells.Item(2,15).Value = 0.3579027849973545
$ws.Cells.Item(2,16).Value = 0.3579027849973545
$ws.Cells.Item(2,17).Value = 37.51549495451511
$ws.Cells.Item(2,18).Value = 337.639454590636
$ws.Cells.Item(2,19).Value = 0.213568798233565
$ws.Cells.Item(2,20).Value = 0.213568798233565

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "L1cam"
$ws.Cells.Item(3,3).Value = "Erbb2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 11.88712433333333
$ws.Cells.Item(3,8).Value = 35.661373
$ws.Cells.Item(3,9).Value = 0.5967229292030898
$ws.Cells.Item(3,10).Value = 0.5967229292030898
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 3.165953666666667
$ws.Cells.Item(3,14).Value = 9.497861
$ws.Cells.Item(3,15).Value = 0.359034148472735
$ws.Cells.Item(3,16).Value = 0.359034148472735
$ws.Cells.Item(3,17).Value = 37.63408486923922
$ws.Cells.Item(3,18).Value = 338.706763823153
$ws.Cells.Item(3,19).Value = 0.2142439087605874
$ws.Cells.Item(3,20).Value = 0.2142439087605874

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "L1cam"
$ws.Cells.Item(4,3).Value = "Erbb2"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 11.88712433333333
$ws.Cells.Item(4,8).Value = 35.661373
$ws.Cells.Item(4,9).Value = 0.5967229292030898
$ws.Cells.Item(4,10).Value = 0.5967229292030898
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.496042666666666
$ws.Cells.Item(4,14).Value = 7.488128
$ws.Cells.Item(4,15).Value = 0.2830630665299106
$ws.Cells.Item(4,16).Value = 0.2830630665299106
$ws.Cells.Item(4,17).Value = 29.67076951997155
$ws.Cells.Item(4,18).Value = 267.0369256797439
$ws.Cells.Item(4,19).Value = 0.1689102222089374
$ws.Cells.Item(4,20).Value = 0.1689102222089374

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "L1cam"
$ws.Cells.Item(5,3).Value = "Erbb2"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.3785896666666667
$ws.Cells.Item(5,8).Value = 1.135769
$ws.Cells.Item(5,9).Value = 0.01900486009268527
$ws.Cells.Item(5,10).Value = 0.01900486009268527
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 3.155977333333333
$ws.Cells.Item(5,14).Value = 9.467931999999999
$ws.Cells.Item(5,15).Value = 0.3579027849973545
$ws.Cells.Item(5,16).Value = 0.3579027849973545
$ws.Cells.Item(5,17).Value = 1.194820406634222
$ws.Cells.Item(5,18).Value = 10.753383659708
$ws.Cells.Item(5,19).Value = 0.00680189235565714
$ws.Cells.Item(5,20).Value = 0.00680189235565714

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "L1cam"
$ws.Cells.Item(6,3).Value = "Erbb2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 0.6666666666666666
$ws.Cells.Item(6,7).Value = 0.3785896666666667
$ws.Cells.Item(6,8).Value = 1.135769
$ws.Cells.Item(6,9).Value = 0.01900486009268527
$ws.Cells.Item(6,10).Value = 0.01900486009268527
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 3.165953666666667
$ws.Cells.Item(6,14).Value = 9.497861
$ws.Cells.Item(6,15).Value = 0.359034148472735
$ws.Cells.Item(6,16).Value = 0.359034148472735
$ws.Cells.Item(6,17).Value = 1.198597343345444
$ws.Cells.Item(6,18).Value = 10.787376090109
$ws.Cells.Item(6,19).Value = 0.00682339376022072
$ws.Cells.Item(6,20).Value = 0.00682339376022072

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "L1cam"
$ws.Cells.Item(7,3).Value = "Erbb2"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 0.6666666666666666
$ws.Cells.Item(7,7).Value = 0.3785896666666667
$ws.Cells.Item(7,8).Value = 1.135769
$ws.Cells.Item(7,9).Value = 0.01900486009268527
$ws.Cells.Item(7,10).Value = 0.01900486009268527
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 2.496042666666666
$ws.Cells.Item(7,14).Value = 7.488128
$ws.Cells.Item(7,15).Value = 0.2830630665299106
$ws.Cells.Item(7,16).Value = 0.2830630665299106
$ws.Cells.Item(7,17).Value = 0.944975961159111
$ws.Cells.Item(7,18).Value = 8.504783650432
$ws.Cells.Item(7,19).Value = 0.005379573976807416
$ws.Cells.Item(7,20).Value = 0.005379573976807416

# Row 8
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "L1cam"
$ws.Cells.Item(8,3).Value = "Erbb2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 4.123197999999999
$ws.Cells.Item(8,8).Value = 12.369594
$ws.Cells.Item(8,9).Value = 0.2069808238940481
$ws.Cells.Item(8,10).Value = 0.2069808238940482
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 3.155977333333333
$ws.Cells.Item(8,14).Value = 9.467931999999999
$ws.Cells.Item(8,15).Value = 0.3579027849973545
$ws.Cells.Item(8,16).Value = 0.3579027849973545
$ws.Cells.Item(8,17).Value = 13.01271942884533
$ws.Cells.Item(8,18).Value = 117.114474859608
$ws.Cells.Item(8,19).Value = 0.0740790133127268
$ws.Cells.Item(8,20).Value = 0.07407901331272682

# Row 9
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "L1cam"
$ws.Cells.Item(9,3).Value = "Erbb2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 4.123197999999999
$ws.Cells.Item(9,8).Value = 12.369594
$ws.Cells.Item(9,9).Value = 0.2069808238940481
$ws.Cells.Item(9,10).Value = 0.2069808238940482
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 3.165953666666667
$ws.Cells.Item(9,14).Value = 9.497861
$ws.Cells.Item(9,15).Value = 0.359034148472735
$ws.Cells.Item(9,16).Value = 0.359034148472735
$ws.Cells.Item(9,17).Value = 13.05385382649266
$ws.Cells.Item(9,18).Value = 117.484684438434
$ws.Cells.Item(9,19).Value = 0.0743131838569847
$ws.Cells.Item(9,20).Value = 0.0743131838569847

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "L1cam"
$ws.Cells.Item(10,3).Value = "Erbb2"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 4.123197999999999
$ws.Cells.Item(10,8).Value = 12.369594
$ws.Cells.Item(10,9).Value = 0.2069808238940481
$ws.Cells.Item(10,10).Value = 0.2069808238940482
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 2.496042666666666
$ws.Cells.Item(10,14).Value = 7.488128
$ws.Cells.Item(10,15).Value = 0.2830630665299106
$ws.Cells.Item(10,16).Value = 0.2830630665299106
$ws.Cells.Item(10,17).Value = 10.29167813111467
$ws.Cells.Item(10,18).Value = 92.62510318003199
$ws.Cells.Item(10,19).Value = 0.05858862672433667
$ws.Cells.Item(10,20).Value = 0.05858862672433668

# Row 11
$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "L1cam"
$ws.Cells.Item(11,3).Value = "Erbb2"
$ws.Cells.Item(11,4).Value = "ECs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 3.531764333333333
$ws.Cells.Item(11,8).Value = 10.595293
$ws.Cells.Item(11,9).Value = 0.1772913868101768
$ws.Cells.Item(11,10).Value = 0.1772913868101768
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 3.155977333333333
$ws.Cells.Item(11,14).Value = 9.467931999999999
$ws.Cells.Item(11,15).Value = 0.3579027849973545
$ws.Cells.Item(11,16).Value = 0.3579027849973545
$ws.Cells.Item(11,17).Value = 11.14616818267511
$ws.Cells.Item(11,18).Value = 100.315513644076
$ws.Cells.Item(11,19).Value = 0.0634530810954055
$ws.Cells.Item(11,20).Value = 0.0634530810954055

# Row 12
$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "L1cam"
$ws.Cells.Item(12,3).Value = "Erbb2"
$ws.Cells.Item(12,4).Value = "FAPs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 3.531764333333333
$ws.Cells.Item(12,8).Value = 10.595293
$ws.Cells.Item(12,9).Value = 0.1772913868101768
$ws.Cells.Item(12,10).Value = 0.1772913868101768
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 3.165953666666667
$ws.Cells.Item(12,14).Value = 9.497861
$ws.Cells.Item(12,15).Value = 0.359034148472735
$ws.Cells.Item(12,16).Value = 0.359034148472735
$ws.Cells.Item(12,17).Value = 11.18140224091922
$ws.Cells.Item(12,18).Value = 100.632620168273
$ws.Cells.Item(12,19).Value = 0.06365366209494208
$ws.Cells.Item(12,20).Value = 0.06365366209494208

# Row 13
$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "L1cam"
$ws.Cells.Item(13,3).Value = "Erbb2"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 3.531764333333333
$ws.Cells.Item(13,8).Value = 10.595293
$ws.Cells.Item(13,9).Value = 0.1772913868101768
$ws.Cells.Item(13,10).Value = 0.1772913868101768
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 2.496042666666666
$ws.Cells.Item(13,14).Value = 7.488128
$ws.Cells.Item(13,15).Value = 0.2830630665299106
$ws.Cells.Item(13,16).Value = 0.2830630665299106
$ws.Cells.Item(13,17).Value = 8.815434464611554
$ws.Cells.Item(13,18).Value = 79.33891018150399
$ws.Cells.Item(13,19).Value = 0.05018464361982918
$ws.Cells.Item(13,20).Value = 0.05018464361982918

